$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source inlineStr formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '67.781.37'
$ws.Range("E2").Value = '  -2.61%  '
$ws.Range("D3").Value = '3.559.61'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '615.69'
$ws.Range("E5").Value = '  -7.38%  '
$ws.Range("D6").Value = '153.89'
$ws.Range("E6").Value = '  -3.84%  '
$ws.Range("D7").Value = '3.557.13'
$ws.Range("E7").Value = '  -3.44%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("E11").Value = '  -3.37%  '
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").Value = '4.164.93'
$ws.Range("D15").Value = '32.04'
$ws.Range("E15").Value = '  -2.19%  '
$ws.Range("D16").Value = '3.553.31'
$ws.Range("E16").Value = '  -3.73%  '
$ws.Range("D17").Value = '67.846.27'
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").Value = '15.64'
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("D21").Value = '454.39'
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("D22").Value = '9.65'
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("D23").Value = '0.645'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '77.64'
$ws.Range("E24").Value = '  -2.78%  '
$ws.Range("D25").Value = '3.704.92'
$ws.Range("E25").Value = '  -3.35%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '10.63'
$ws.Range("E27").Value = '  -2.85%  '
$ws.Range("E28").Value = '  -8.14%  '
$ws.Range("D29").Value = '8.38'
$ws.Range("E29").Value = '  -7.18%  '
$ws.Range("E30").Value = '  -3.86%  '
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = '25.94'
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  -4.71%  '
$ws.Range("D35").Value = '6.21'
$ws.Range("E35").Value = '  -4.02%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.561.76'
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.157'
$ws.Range("E37").Value = '  -4.27%  '
$ws.Range("D38").Value = '8.07'
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '176.60'
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").Value = '0.0884'
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("D43").Value = '5.61'
$ws.Range("E43").Value = '  -7.86%  '
$ws.Range("D44").Value = '2.09'
$ws.Range("E44").Value = '  -6.82%  '
$ws.Range("D45").Value = '0.893'
$ws.Range("E45").Value = '  -4.29%  '
$ws.Range("D46").Value = '46.25'
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").Value = '28.66'
$ws.Range("E47").Value = '  +4.13%  '
$ws.Range("E48").Value = '  -5.84%  '
$ws.Range("D49").Value = '7.70'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").Value = '1.20'
$ws.Range("E50").Value = '  -6.42%  '
$ws.Range("E51").Value = '  -5.35%  '
